$wb = $excel.ActiveWorkbook

# --- Rename the two "25 GW - High (..)" sheets (formulas that self-reference
# their own sheet name, e.g. '25 GW - High (SC)'!M2, will auto-update) ---
$wsSC = $wb.Worksheets.Item("25 GW - High (SC)")
$wsSC.Name = "25 GW (SC)"

$wsCC = $wb.Worksheets.Item("25 GW - High (CC)")
$wsCC.Name = "25 GW (CC)"

# --- Update per-sheet selections ---
# "25 GW (CC)" (formerly "25 GW - High (CC)"): selection moves from E11 to H4,
# and this sheet is NOT the active tab in the final state.
$wsCC.Activate()
$wsCC.Range("H4").Select()

# "25 GW (SC)" (formerly "25 GW - High (SC)"): selection moves from C26 to G19,
# and this sheet becomes the active tab (tabSelected), replacing "55 GW".
$wsSC.Activate()
$wsSC.Range("G19").Select()
